{"js": "// Hybrid bold + color highlighting for quantitative impact metrics\n// (percentages, dollar amounts, large numbers) in the resume body.\n//\n// For each targeted paragraph we locate it by its exact original text,\n// then split out the specific numeric substrings into their own runs\n// with bold + color (#2C3E50) formatting, leaving the surrounding text\n// untouched (same run formatting it already had).\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Each entry: the paragraph's ORIGINAL full text (used to find the right\n// paragraph unambiguously) and the ordered list of substrings inside it\n// that must become bold + colored. Every substring occurs exactly once\n// within its paragraph, so a simple in-order search works.\nconst targets = [\n  {\n    paragraphText:\n      \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    highlights: [\"23%\", \"64%\"],\n  },\n  {\n    paragraphText:\n      \"\\u2022 Utilized advanced sampling methods to decrease survey margin of error from \\u00B14.2% to \\u00B12.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes\",\n    highlights: [\"\\u00B14.2%\", \"\\u00B12.1%\", \"71%\", \"87%\"],\n  },\n  {\n    paragraphText:\n      \"\\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n    highlights: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    paragraphText:\n      \"\\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion\",\n    highlights: [\"$2\"],\n  },\n  {\n    paragraphText:\n      \"\\u2022 Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%\",\n    highlights: [\"57%\"],\n  },\n  {\n    paragraphText:\n      \"\\u2022 Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from \\u00B14.2% to \\u00B12.1%\",\n    highlights: [\"\\u00B14.2%\", \"\\u00B12.1%\"],\n  },\n  {\n    paragraphText: \"\\u2022 Increased voter turnout prediction accuracy from 71% to 87%\",\n    highlights: [\"71%\", \"87%\"],\n  },\n  {\n    paragraphText:\n      \"\\u2022 Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%\",\n    highlights: [\"34%\", \"28%\"],\n  },\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const target of targets) {\n  const para = paragraphs.items.find((p) => p.text === target.paragraphText);\n  if (!para) {\n    // Should always be found; skip defensively if the doc doesn't match.\n    continue;\n  }\n  for (const needle of target.highlights) {\n    const found = para.search(needle, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n    for (let i = 0; i < found.items.length; i++) {\n      found.items[i].font.bold = true;\n      found.items[i].font.color = HIGHLIGHT_COLOR;\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Hybrid bold + color highlighting for quantitative impact metrics\n# (percentages, dollar amounts, large numbers) in the resume body.\n#\n# For each targeted paragraph (matched by its exact original text) we\n# walk the list of numeric substrings that must become bold + colored\n# (#2C3E50 == OLE/COM color 5258796, since Word COM colors are stored\n# BGR: 0x50 3E 2C -> R=2C G=3E B=50) and apply character formatting to\n# just that substring, leaving the rest of the paragraph's run(s) as-is.\n\n$d = $word.ActiveDocument\n$bullet = [char]0x2022\n$plusMinus = [char]0x00B1\n$highlightColor = 5258796\n\n$targets = @(\n    @{\n        text = \"$bullet Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n        highlights = @(\"23%\", \"64%\")\n    },\n    @{\n        text = \"$bullet Utilized advanced sampling methods to decrease survey margin of error from ${plusMinus}4.2% to ${plusMinus}2.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes\"\n        highlights = @(\"${plusMinus}4.2%\", \"${plusMinus}2.1%\", \"71%\", \"87%\")\n    },\n    @{\n        text = \"$bullet Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis\"\n        highlights = @(\"73.5%\", \"`$4.7M\")\n    },\n    @{\n        text = \"$bullet Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion\"\n        highlights = @(\"`$2\")\n    },\n    @{\n        text = \"$bullet Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%\"\n        highlights = @(\"57%\")\n    },\n    @{\n        text = \"$bullet Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from ${plusMinus}4.2% to ${plusMinus}2.1%\"\n        highlights = @(\"${plusMinus}4.2%\", \"${plusMinus}2.1%\")\n    },\n    @{\n        text = \"$bullet Increased voter turnout prediction accuracy from 71% to 87%\"\n        highlights = @(\"71%\", \"87%\")\n    },\n    @{\n        text = \"$bullet Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%\"\n        highlights = @(\"34%\", \"28%\")\n    }\n)\n\nforeach ($target in $targets) {\n    foreach ($p in $d.Paragraphs) {\n        $ptext = $p.Range.Text.TrimEnd([char]0x0D)\n        if ($ptext -eq $target.text) {\n            $paraEnd = $p.Range.End\n            $cursor = $p.Range.Start\n            foreach ($needle in $target.highlights) {\n                $r = $d.Range($cursor, $paraEnd)\n                $find = $r.Find\n                $find.Text = $needle\n                $find.MatchCase = $true\n                $found = $find.Execute()\n                if ($found) {\n                    $r.Bold = 1\n                    $r.Font.Color = $highlightColor\n                    $cursor = $r.End\n                }\n            }\n            break\n        }\n    }\n}\n"}
